$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.086.62"
$ws.Range("E2").Value = "  -0.17%  "

$ws.Range("D3").Value = "1.898.39"
$ws.Range("E3").Value = "  -0.42%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.22%  "

$ws.Range("D5").Value = "325.36"
$ws.Range("E5").Value = "  -0.44%  "

$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").Value = "  -0.23%  "

$ws.Range("D7").Value = "0.4604"
$ws.Range("E7").Value = "  -0.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07860"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.95%  "

$ws.Range("D10").Value = "0.9907"
$ws.Range("E10").Value = "  -0.88%  "

$ws.Range("D11").Value = "21.85"
$ws.Range("E11").Value = "  -1.64%  "

$ws.Range("D12").Value = "1.870.99"
$ws.Range("E12").Value = "  -1.89%  "

$ws.Range("D13").Value = "5.772"
$ws.Range("E13").Value = "  +0.19%  "

$ws.Range("E14").Value = "  -0.56%  "

$ws.Range("D15").Value = "0.06988"
$ws.Range("E15").Value = "  +0.79%  "

$ws.Range("D16").Value = "87.86"
$ws.Range("E16").Value = "  -0.57%  "

$ws.Range("E17").Value = "  -0.21%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009930"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.32%  "

$ws.Range("E19").Value = "  -0.48%  "

$ws.Range("D20").Value = "0.9989"
$ws.Range("E20").Value = "  -0.31%  "

$ws.Range("D21").Value = "29.085.64"
$ws.Range("E21").Value = "  -0.27%  "

$ws.Range("E22").Value = "  -0.79%  "

$ws.Range("E23").Value = "  -0.01%  "

$ws.Range("D24").Value = "2.100.68"
$ws.Range("E24").Value = "  -1.71%  "

$ws.Range("E25").Value = "  +1.88%  "

$ws.Range("D26").Value = "155.85"
$ws.Range("E26").Value = "  -0.46%  "

$ws.Range("D27").Value = "19.43"
$ws.Range("E27").Value = "  +0.02%  "

$ws.Range("D28").Value = "5.878"
$ws.Range("E28").Value = "  -3.64%  "

$ws.Range("D29").Value = "118.54"
$ws.Range("E29").Value = "  -0.37%  "

$ws.Range("E30").Value = "  -5.90%  "

$ws.Range("D31").Value = "0.09313"
$ws.Range("E31").Value = "  -0.71%  "

$ws.Range("D32").Value = "0.8994"
$ws.Range("E32").Value = "  -2.26%  "

$ws.Range("D33").Value = "5.236"
$ws.Range("E33").Value = "  -1.73%  "

$ws.Range("D34").Value = "1.322"
$ws.Range("E34").Value = "  -2.02%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.150"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.53%  "

$ws.Range("D36").Value = "0.05798"
$ws.Range("E36").Value = "  -0.31%  "

$ws.Range("D37").Value = "1.171"
$ws.Range("E37").Value = "  -2.71%  "

$ws.Range("D38").Value = "0.02081"
$ws.Range("E38").Value = "  -1.04%  "

$ws.Range("D39").Value = "0.9989"
$ws.Range("E39").Value = "  -0.28%  "

$ws.Range("D40").Value = "7.679"
$ws.Range("E40").Value = "  -3.26%  "

$ws.Range("D41").Value = "0.5674"
$ws.Range("E41").Value = "  -1.45%  "

$ws.Range("D42").Value = "0.1793"
$ws.Range("E42").Value = "  -0.04%  "

$ws.Range("D43").Value = "9.724"
$ws.Range("E43").Value = "  -2.21%  "

$ws.Range("D44").Value = "11.92"
$ws.Range("E44").Value = "  -0.81%  "

$ws.Range("D45").Value = "2.235"
$ws.Range("E45").Value = "  +0.89%  "

$ws.Range("D46").Value = "0.5335"
$ws.Range("E46").Value = "  -1.65%  "

$ws.Range("E47").Value = "  -1.04%  "

$ws.Range("D48").Value = "1.849"
$ws.Range("E48").Value = "  -1.24%  "

$ws.Range("D49").Value = "2.551"
$ws.Range("E49").Value = "  +0.08%  "

$ws.Range("D50").Value = "112.82"
$ws.Range("E50").Value = "  +0.59%  "

$ws.Range("D51").Value = "1.044"
$ws.Range("E51").Value = "  -2.10%  "

